$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicators")

$ws.Range("B5").Value = 0.86486486486486491
$ws.Range("B6").Value = 0.82051282051282048
$ws.Range("B7").Value = 0.84210526315789469

$excel.ActiveWindow.WindowState = -4143
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 16800
$excel.ActiveWindow.Height = 15600
